# Adds a new "#TUS" column of data to the "Tasks1&2" sheet.
#
# The new column is inserted between the existing "QLTY" (G) and "TEST" (H)
# columns, which means every column from the old H ("TEST") through the old
# U ("Q3") needs to shift one position to the right (H->I, I->J, ... U->V),
# and the new column H is populated with the "#TUS" header/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift columns H(8)..U(21) one column to the right, for the
# header row and the 18 data rows (rows 1-19). Rows 20-22 only contain a
# handful of blank styled filler cells that never reach column H in a way
# that needs shifting, so they are intentionally left untouched.
# Process right-to-left (U first, H last) so that each destination column is
# written before its own original content is overwritten by the next step.
for ($colIndex = 21; $colIndex -ge 8; $colIndex--) {
    $srcRange = $ws.Range($ws.Cells.Item(1, $colIndex), $ws.Cells.Item(19, $colIndex))
    $destTopLeft = $ws.Cells.Item(1, $colIndex + 1)
    $srcRange.Copy($destTopLeft)
}

# --- Step 2: the copy above turns formulas into static values, so restore
# the two shared formulas (now living in columns M and P instead of the
# original L and O) that sum the shifted PROD/TEST pairs.
$ws.Range("M2:M19").Formula = "=SUM(K2:L2)"
$ws.Range("P2:P19").Formula = "=SUM(N2:O2)"

# --- Step 3: populate the freshly emptied column H with the new "#TUS"
# header and its per-row values.
$ws.Range("H1").Value = "#TUS"

$tusValues = @(5, 5, 5, 3, 5, 3, 5, 2, 5, 1, 5, 0, 5, 5, 5, 2, 5, 3)
for ($i = 0; $i -lt $tusValues.Length; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 8).Value = $tusValues[$i]
}

# --- Step 4: match the final selected cell recorded in the sheet view.
$ws.Range("I22").Select()
